# Update countries & provincias Spain
# Applies the 9-Abril-2020 17:52 data refresh on top of the 17:22 snapshot:
#  - reorders "Cuba" ahead of "Costa Rica" (so row 87 becomes Cuba, 88 Costa
#    Rica, 89 Afganistan) and refreshes the numbers for the affected rows
#  - refreshes case/death counters for a handful of other countries
#  - bumps the "Datos actualizados ..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 17:52"

# --- plain numeric refreshes (country stays in the same row) ----------
$ws.Range("B4").Value  = 436969
$ws.Range("C4").Value  = 2042
$ws.Range("D4").Value  = 24391
$ws.Range("E4").Value  = 396870
$ws.Range("G4").Value  = 920
$ws.Range("H4").Value  = 15708

$ws.Range("B16").Value = 19805
$ws.Range("C16").Value = 367
$ws.Range("D16").Value = 4884
$ws.Range("E16").Value = 14459
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = 462

$ws.Range("B17").Value = 16474
$ws.Range("C17").Value = 286
$ws.Range("E17").Value = 15462
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 839

$ws.Range("B31").Value = 5467
$ws.Range("C31").Value = 155
$ws.Range("D31").Value = 301
$ws.Range("E31").Value = 5054
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 112

$ws.Range("E32").Value = 4309
$ws.Range("G32").Value = 26
$ws.Range("H32").Value = 246

$ws.Range("B34").Value = 4489
$ws.Range("C34").Value = 226
$ws.Range("E34").Value = 3854

$ws.Range("B42").Value = 3115
$ws.Range("C42").Value = 81
$ws.Range("E42").Value = 2563
$ws.Range("F42").Value = 30
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 52

$ws.Range("D55").Value = 365
$ws.Range("E55").Value = 1360

$ws.Range("B56").Value = 1648
$ws.Range("C56").Value = 32
$ws.Range("D56").Value = 688
$ws.Range("E56").Value = 954
$ws.Range("F56").Value = 11

$ws.Range("B74").Value = 846
$ws.Range("C74").Value = 42
$ws.Range("E74").Value = 710

$ws.Range("E82").Value = 570
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 3

$ws.Range("B112").Value = 218
$ws.Range("C112").Value = 7
$ws.Range("E112").Value = 165

# --- Republica de Chipre (row 86) keeps its place, numbers refresh ----
$ws.Range("B86").Value = 564
$ws.Range("C86").Value = 38
$ws.Range("D86").Value = 53
$ws.Range("E86").Value = 501
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 10

# --- Cuba moves up, ahead of Costa Rica/Afganistan, with fresh numbers-
$ws.Range("A87").Value = "Cuba"
$ws.Range("B87").Value = 515
$ws.Range("C87").Value = 58
$ws.Range("D87").Value = 28
$ws.Range("E87").Value = 472
$ws.Range("F87").Value = 15
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 15

# Costa Rica drops to row 88, keeping its existing (unchanged) figures
$ws.Range("A88").Value = "Costa Rica"
$ws.Range("B88").Value = 502
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 29
$ws.Range("E88").Value = 470
$ws.Range("F88").Value = 15
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 3

# Afganistan drops to row 89, keeping its existing (unchanged) figures
$ws.Range("A89").Value = "Afganistan"
$ws.Range("B89").Value = 484
$ws.Range("C89").Value = 40
$ws.Range("D89").Value = 32
$ws.Range("E89").Value = 437
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 15
